$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "#" header cell and the "Massstab"/"RS" column that is no
# longer used (Schueler(in) now spans B:C on each of these rows).
$ws.Range("A1").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("C14").Value = ""

# Merge B:C on the rows that previously held the now-removed "Massstab"/"RS"
# values, so "Schueler(in)" / student names span the merged cell.
$ws.Range("B8:C8").Merge()
$ws.Range("B9:C9").Merge()
$ws.Range("B11:C11").Merge()
$ws.Range("B12:C12").Merge()
$ws.Range("B13:C13").Merge()
$ws.Range("B14:C14").Merge()

# The data cells in the border=2 style now center their (empty) content.
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("C8").VerticalAlignment = -4108

# Restore the selection to A1 (it had drifted to C4).
$ws.Range("A1").Select()
